# Update cryptos list - prices (col D) and volume/1h changes (col E)
# Row 43/44 also swap coin identity (VeChain <-> InjectiveProtocol)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.061.44"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.92%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.999.31"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.03%  "
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.642"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.21%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.22"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +13.17%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "58.38"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.55%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.368"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.47%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0746"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.05%  "
$ws.Range("E12").Value = "  -1.31%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.947"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.87"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.34%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.284.29"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.39%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.42"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.36%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.51"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +15.99%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.983.40"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.76%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "35.957.38"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.82%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.01"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.47%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0853"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.70%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.25"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.67%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "233.71"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.66%  "
$ws.Range("E24").Value = "  +0.19%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.63"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +18.10%  "
$ws.Range("E26").Value = "  -3.38%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.70"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +7.45%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "165.74"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.33%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.61"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.35%  "
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.10"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.30%  "
$ws.Range("E32").Value = "  -0.97%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0996"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +14.88%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0606"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.62%  "
$ws.Range("E35").Value = "  +2.74%  "
$ws.Range("E36").Value = "  +12.24%  "
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.79"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.50%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.78"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +17.41%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.24"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.56%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0964"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +8.36%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.89"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.30%  "
$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.96"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +9.57%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0215"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.81%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.12"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.81%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "94.17"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.79"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.13%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.367.60"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.30%  "
$ws.Range("E49").Value = "  -0.34%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.34"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.40%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "46.76"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.04%  "
